$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("testCitizen")

# Update the country-code style labels in column A (rows 1-8) of testCitizen
$ws2.Range("A1").Value = "ulkemb3is11"
$ws2.Range("A2").Value = "ulkemb3is12"
$ws2.Range("A3").Value = "ulkemb3is13"
$ws2.Range("A4").Value = "ulkemb3is14"
$ws2.Range("A5").Value = "ulkemb3is15"
$ws2.Range("A6").Value = "ulkemb3is16"
$ws2.Range("A7").Value = "ulkemb3is17"
$ws2.Range("A8").Value = "ulkemb3is18"

# Update the codes in column B (rows 1-8) of testCitizen
$ws2.Range("B1").Value = "ub3is11"
$ws2.Range("B2").Value = "ub3is12"
$ws2.Range("B3").Value = "ub3is13"
$ws2.Range("B4").Value = "ub3is14"
$ws2.Range("B5").Value = "ub3is15"
$ws2.Range("B6").Value = "ub3is16"
$ws2.Range("B7").Value = "ub3is17"
$ws2.Range("B8").Value = "ub3is18"

# Update the active selection on testCitizen sheet
$ws2.Activate()
$ws2.Range("B1:B8").Select()
